$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Delete the "Meta description: ..." paragraph that currently sits
#    right under the "Play Bonsai Spins Free Slot Machine Online"
#    H1 heading.
# ---------------------------------------------------------------------
$metaPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.StartsWith("Meta description")) {
        $metaPara = $para
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------
# 2) Find the closing "Prompt: Create a feature image..." paragraph and:
#      a) insert a new bold "Play Bonsai Spins Free Slot Machine Online"
#         paragraph right before it
#      b) replace its own text with the meta-description sentence,
#         keeping the existing italic run formatting
# ---------------------------------------------------------------------
$promptPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.StartsWith("Prompt: Create a feature image")) {
        $promptPara = $para
        break
    }
}

if ($promptPara -ne $null) {
    $insertRange = $d.Range($promptPara.Range.Start, $promptPara.Range.End - 1)
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Bonsai Spins Free Slot Machine Online</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Discover the exciting Bonsai Spins slot machine game with Wild and Scatter symbols. Play now for free and experience beautiful graphics and Oriental music.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insertRange.InsertXML($xml)
}
